$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(216).EntireRow.Insert(-4121)
Write-Host "done"
